# Recomputed "2023" metrics (Drop2023/Wday2023/Wend2023/High2023/Low2023 -- columns F,J,N,R,V)
# for every station row, plus a full-row correction for row 10 (ReneLevesque_Wolfe).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 36.2
$ws.Range("J2").Value = 29.85
$ws.Range("N2").Value = 59.9
$ws.Range("R2").Value = 31.91
$ws.Range("V2").Value = 121.07

# Row 3
$ws.Range("F3").Value = 34.14
$ws.Range("J3").Value = 27.87
$ws.Range("N3").Value = 66.24
$ws.Range("R3").Value = 29.11
$ws.Range("V3").Value = 118.77

# Row 4
$ws.Range("F4").Value = -13.14
$ws.Range("J4").Value = -15.42
$ws.Range("N4").Value = -4.14
$ws.Range("R4").Value = -16.18
$ws.Range("V4").Value = 46.56

# Row 5
$ws.Range("F5").Value = -98.25
$ws.Range("J5").Value = -98.14
$ws.Range("N5").Value = -98.75
$ws.Range("R5").Value = -100
$ws.Range("V5").Value = -72.4

# Row 6
$ws.Range("F6").Value = 4.03
$ws.Range("J6").Value = -1.35
$ws.Range("N6").Value = 25.87
$ws.Range("R6").Value = 1.08
$ws.Range("V6").Value = 86.39

# Row 7
$ws.Range("F7").Value = 16.62
$ws.Range("J7").Value = 9.86
$ws.Range("N7").Value = 57.09
$ws.Range("R7").Value = 13.39
$ws.Range("V7").Value = 68.7

# Row 8
$ws.Range("F8").Value = 11.76
$ws.Range("J8").Value = 10.71
$ws.Range("N8").Value = 14.55
$ws.Range("R8").Value = 8.83
$ws.Range("V8").Value = 62.88

# Row 9
$ws.Range("F9").Value = -11.52
$ws.Range("J9").Value = -15.68
$ws.Range("N9").Value = 2.22
$ws.Range("R9").Value = -13.35
$ws.Range("V9").Value = 13.14

# Row 10
$ws.Range("C10").Value = -23.95
$ws.Range("D10").Value = -9.07
$ws.Range("E10").Value = -19.02
$ws.Range("F10").Value = -21.31
$ws.Range("G10").Value = -29.37
$ws.Range("H10").Value = -13.89
$ws.Range("I10").Value = -24.16
$ws.Range("J10").Value = -25.02
$ws.Range("K10").Value = -4.65
$ws.Range("L10").Value = 8.07
$ws.Range("M10").Value = -0.77
$ws.Range("N10").Value = -8.14
$ws.Range("O10").Value = -26.7
$ws.Range("P10").Value = -8.98
$ws.Range("Q10").Value = -18.54
$ws.Range("R10").Value = -21.71
$ws.Range("S10").Value = 16.45
$ws.Range("T10").Value = -10.31
$ws.Range("U10").Value = -26.2
$ws.Range("V10").Value = -15.44

# Row 11
$ws.Range("F11").Value = -26.1
$ws.Range("J11").Value = -24.41
$ws.Range("N11").Value = -32.88
$ws.Range("R11").Value = -29.4
$ws.Range("V11").Value = 24.46

# Row 12
$ws.Range("F12").Value = -30.53
$ws.Range("J12").Value = -41.47
$ws.Range("N12").Value = 34.31
$ws.Range("R12").Value = -35.22
$ws.Range("V12").Value = 46.62

# Row 14
$ws.Range("F14").Value = -13.59
$ws.Range("J14").Value = -21.23
$ws.Range("N14").Value = 2.48
$ws.Range("R14").Value = -14.57
$ws.Range("V14").Value = 100.15

# Row 15
$ws.Range("F15").Value = -91.02
$ws.Range("J15").Value = -90.5
$ws.Range("N15").Value = -92.46
$ws.Range("R15").Value = -94.25
$ws.Range("V15").Value = 41.86

# Row 16
$ws.Range("F16").Value = -33.52
$ws.Range("J16").Value = -35.96
$ws.Range("N16").Value = -25.78
$ws.Range("R16").Value = -35.22
$ws.Range("V16").Value = -0.89

# Row 17
$ws.Range("F17").Value = -64.3
$ws.Range("J17").Value = -66.33
$ws.Range("N17").Value = -58.01
$ws.Range("R17").Value = -65.55
$ws.Range("V17").Value = -39.24

# Row 18
$ws.Range("F18").Value = -59.09
$ws.Range("J18").Value = -57.7
$ws.Range("N18").Value = -62.67
$ws.Range("R18").Value = -59.47
$ws.Range("V18").Value = -31.36

# Row 19
$ws.Range("F19").Value = -8.18
$ws.Range("J19").Value = -13.81
$ws.Range("N19").Value = 8.99
$ws.Range("R19").Value = -10.96
$ws.Range("V19").Value = 91.25

# Row 20
$ws.Range("F20").Value = -12.48
$ws.Range("J20").Value = -15.05
$ws.Range("N20").Value = -1.98
$ws.Range("R20").Value = -15.5
$ws.Range("V20").Value = 37.07

# Row 21
$ws.Range("F21").Value = -42.41
$ws.Range("J21").Value = -45.22
$ws.Range("N21").Value = -31.24
$ws.Range("R21").Value = -43.83
$ws.Range("V21").Value = -18.48
